# Auto-generated edit script: applies scheduled market-data refresh values
# to the Ifrit_Profits crafting-leve profit sheets (one worksheet per job).
$wb = $excel.ActiveWorkbook

# ---- ALC sheet ----
$ws = $wb.Worksheets.Item("ALC")
# Row 18
$ws.Range("H18").Value = 990.1
$ws.Range("I18").Value = 990.1
$ws.Range("K18").Value = 990.1
$ws.Range("M18").Value = -706.1
# Row 100
$ws.Range("H100").Value = 1628.1364
$ws.Range("I100").Value = 1454.6
$ws.Range("K100").Value = 1454.6
$ws.Range("M100").Value = -913.5999999999999

# ---- ARM sheet ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 10503
$ws.Range("I32").Value = 6959.909
$ws.Range("K32").Value = 6959.909
$ws.Range("M32").Value = -6672.909
# Row 110
$ws.Range("H110").Value = 1151.7894
$ws.Range("I110").Value = 852.1539
$ws.Range("J110").Value = 1801
$ws.Range("K110").Value = 852.1539
$ws.Range("L110").Value = 1801
$ws.Range("M110").Value = 1192.8461
$ws.Range("N110").Value = -5891

# ---- BSM sheet ----
$ws = $wb.Worksheets.Item("BSM")
# Row 64
$ws.Range("H64").Value = 11000
$ws.Range("I64").Value = 20000
$ws.Range("J64").Value = 2000
$ws.Range("K64").Value = 20000
$ws.Range("L64").Value = 2000
$ws.Range("M64").Value = -19775
$ws.Range("N64").Value = -2450
# Row 67
$ws.Range("H67").Value = 11000
$ws.Range("I67").Value = 20000
$ws.Range("J67").Value = 2000
$ws.Range("K67").Value = 20000
$ws.Range("L67").Value = 2000
$ws.Range("M67").Value = -19220
$ws.Range("N67").Value = -3560
# Row 86
$ws.Range("H86").Value = 2110.9333
$ws.Range("I86").Value = 1364.591
$ws.Range("J86").Value = 4163.375
$ws.Range("K86").Value = 1364.591
$ws.Range("L86").Value = 4163.375
$ws.Range("M86").Value = -241.5909999999999
$ws.Range("N86").Value = -6409.375
# Row 89
$ws.Range("H89").Value = 2110.9333
$ws.Range("I89").Value = 1364.591
$ws.Range("J89").Value = 4163.375
$ws.Range("K89").Value = 6822.955
$ws.Range("L89").Value = 20816.875
$ws.Range("M89").Value = -1206.955
$ws.Range("N89").Value = -32048.875
# Row 105
$ws.Range("H105").Value = 1454.75
$ws.Range("I105").Value = 1431.6666
$ws.Range("J105").Value = 1477.8334
$ws.Range("K105").Value = 1431.6666
$ws.Range("L105").Value = 1477.8334
$ws.Range("M105").Value = 315.3334
$ws.Range("N105").Value = -4971.8334

# ---- CUL sheet ----
$ws = $wb.Worksheets.Item("CUL")
# Row 45
$ws.Range("H45").Value = 890.5
$ws.Range("I45").Value = 1150
$ws.Range("J45").Value = 838.6
$ws.Range("K45").Value = 3450
$ws.Range("L45").Value = 2515.8
$ws.Range("M45").Value = -2918
$ws.Range("N45").Value = -3579.8
# Row 75
$ws.Range("H75").Value = 3500
$ws.Range("I75").Value = 500
$ws.Range("J75").Value = 4100
$ws.Range("K75").Value = 1500
$ws.Range("L75").Value = 12300
$ws.Range("M75").Value = -502
$ws.Range("N75").Value = -14296
# Row 78
$ws.Range("H78").Value = 3500
$ws.Range("I78").Value = 500
$ws.Range("J78").Value = 4100
$ws.Range("K78").Value = 4500
$ws.Range("L78").Value = 36900
$ws.Range("M78").Value = 492
$ws.Range("N78").Value = -46884
# Row 109
$ws.Range("H109").Value = 2036.25
$ws.Range("I109").Value = 1085.7142
$ws.Range("J109").Value = 2775.5557
$ws.Range("K109").Value = 3257.1426
$ws.Range("L109").Value = 8326.667099999999
$ws.Range("M109").Value = -2217.1426
$ws.Range("N109").Value = -10406.6671
# Row 114
$ws.Range("H114").Value = 1564.35
$ws.Range("I114").Value = 361.77777
$ws.Range("J114").Value = 2548.2727
$ws.Range("K114").Value = 1085.33331
$ws.Range("L114").Value = 7644.8181
$ws.Range("M114").Value = 2168.66669
$ws.Range("N114").Value = -14152.8181
# Row 131
$ws.Range("H131").Value = 1893.6428
$ws.Range("I131").Value = 5504.9
$ws.Range("J131").Value = 1483.2727
$ws.Range("K131").Value = 16514.7
$ws.Range("L131").Value = 4449.8181
$ws.Range("M131").Value = -11474.7
$ws.Range("N131").Value = -14529.8181
# Row 137
$ws.Range("H137").Value = 39143.066
$ws.Range("I137").Value = 4286.6665
$ws.Range("J137").Value = 43016
$ws.Range("K137").Value = 12859.9995
$ws.Range("L137").Value = 129048
$ws.Range("M137").Value = -7759.999500000002
$ws.Range("N137").Value = -139248

# ---- GSM sheet ----
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 5027.5557
$ws.Range("I70").Value = 5049.6
$ws.Range("J70").Value = 5000
$ws.Range("K70").Value = 5049.6
$ws.Range("L70").Value = 5000
$ws.Range("M70").Value = -4779.6
$ws.Range("N70").Value = -5540
# Row 73
$ws.Range("H73").Value = 5027.5557
$ws.Range("I73").Value = 5049.6
$ws.Range("J73").Value = 5000
$ws.Range("K73").Value = 5049.6
$ws.Range("L73").Value = 5000
$ws.Range("M73").Value = -4113.6
$ws.Range("N73").Value = -6872
# Row 80
$ws.Range("H80").Value = 2493.7334
$ws.Range("I80").Value = 2675.125
$ws.Range("J80").Value = 2286.4285
$ws.Range("K80").Value = 2675.125
$ws.Range("L80").Value = 2286.4285
$ws.Range("M80").Value = -1677.125
$ws.Range("N80").Value = -4282.4285
# Row 83
$ws.Range("H83").Value = 2493.7334
$ws.Range("I83").Value = 2675.125
$ws.Range("J83").Value = 2286.4285
$ws.Range("K83").Value = 13375.625
$ws.Range("L83").Value = 11432.1425
$ws.Range("M83").Value = -8383.625
$ws.Range("N83").Value = -21416.1425
# Row 97
$ws.Range("H97").Value = 2129.2856
$ws.Range("I97").Value = 2134.1667
$ws.Range("J97").Value = 2100
$ws.Range("K97").Value = 2134.1667
$ws.Range("L97").Value = 2100
$ws.Range("M97").Value = -1638.1667
$ws.Range("N97").Value = -3092

# ---- WVR sheet ----
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 4994.5713
$ws.Range("I81").Value = 3000
$ws.Range("J81").Value = 5792.4
$ws.Range("K81").Value = 6000
$ws.Range("L81").Value = 11584.8
$ws.Range("M81").Value = -4939
$ws.Range("N81").Value = -13706.8
# Row 84
$ws.Range("H84").Value = 4994.5713
$ws.Range("I84").Value = 3000
$ws.Range("J84").Value = 5792.4
$ws.Range("K84").Value = 30000
$ws.Range("L84").Value = 57924
$ws.Range("M84").Value = -24696
$ws.Range("N84").Value = -68532
# Row 107
$ws.Range("H107").Value = 269.1579
$ws.Range("I107").Value = 211
$ws.Range("J107").Value = 349.125
$ws.Range("K107").Value = 633
$ws.Range("L107").Value = 1047.375
$ws.Range("M107").Value = 1287
$ws.Range("N107").Value = -4887.375
